# "Finished testing, updated the test outputs and the test case excel
# sheet (Test Cases - Login.xlsx)"
#
# All test cases (rows 2-9) were re-run and their Fail/Pass outcome in
# column G flips from "Fail" to "Pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2:G9").Value = "Pass"

# Leave the selection where the author left it when they saved the file.
$ws.Range("G2:G9").Select()
